# Rename the sheet (Sheet1 -> GmailKeywords)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "GmailKeywords"

# ---- Row 1 : header row -------------------------------------------------
$ws.Range("A1").Value = "TestCase Name"
$ws.Range("B1").Value = "#"
$ws.Range("C1").Value = "Type"
$ws.Range("D1").Value = "Keyword"
$ws.Range("E1").Value = "Operation"
$ws.Range("F1").Value = "Data Set"
$ws.Range("G1").Value = "Description"
$ws.Range("H1").Value = "Result"

# ---- Row 2 ---------------------------------------------------------------
$ws.Range("A2").Value = "TC_01"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "Function"
$ws.Range("D2").Value = "Run app"

# ---- Row 3 -----------------------------------------------------------
$ws.Range("B3").Value = 2
$ws.Range("D3").Value = "Login"
$ws.Range("E3").Value = "Set value"
$ws.Range("G3").Value = "Login to gmail site"

# ---- Row 4 -----------------------------------------------------------
$ws.Range("B4").Value = 3
$ws.Range("D4").Value = "Emails"
$ws.Range("E4").Value = "Set value"
$ws.Range("G4").Value = "Send email message"

# ---- Row 5 -----------------------------------------------------------
$ws.Range("B5").Value = 4
$ws.Range("D5").Value = "Delete"
$ws.Range("E5").Value = "Click"
$ws.Range("G5").Value = "Delete email message"

# ---- Row 6 -----------------------------------------------------------
$ws.Range("D6").Value = "Compose"
$ws.Range("E6").Value = "Click"
$ws.Range("G6").Value = "Compose email message"

# ---- Styling ---------------------------------------------------------
# Style 1: left + vertical-center alignment, applied to A2:B2
$ws.Range("A2:B2").HorizontalAlignment = -4131
$ws.Range("A2:B2").VerticalAlignment = -4108

# Style 2: left alignment only, applied to A3:B6
$ws.Range("A3:B6").HorizontalAlignment = -4131

# ---- Column widths -----------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 6.16
$ws.Columns.Item(3).ColumnWidth = 15.6598
$ws.Columns.Item(4).ColumnWidth = 12.6598
$ws.Columns.Item(5).ColumnWidth = 26.3251
$ws.Columns.Item(6).ColumnWidth = 30.8254
$ws.Columns.Item(7).ColumnWidth = 19.8254

# ---- Selection -----------------------------------------------------------
$ws.Range("G8").Select()
